# The edit rotates the content of 8 consecutive "content" paragraphs:
# each paragraph's text becomes the text that used to sit in the *next*
# content paragraph (document order), with the last one wrapping to the
# first. Doing this as a single top-to-bottom pass of Find/Replace would
# make later lookups match text a previous step just inserted, so we
# first swap every source text for a unique placeholder token, then swap
# every placeholder for its real destination text.

$d = $word.ActiveDocument

function Replace-Text($needle, $replacement, [bool]$wild) {
    $d.Content.Find.Execute(
        $needle,        # FindText
        $true,          # MatchCase
        $true,          # MatchWholeWord
        $wild,          # MatchWildcards
        $false,         # MatchSoundsLike
        $false,         # MatchAllWordForms
        $true,          # Forward
        1,              # Wrap (wdFindContinue)
        $false,         # Format
        $replacement,   # ReplaceWith
        2               # Replace (wdReplaceAll)
    ) | Out-Null
}

# ---- phase 1: park every source paragraph's text behind a unique token ----
Replace-Text "Propiciar uma integração entre os elementos de estruturação da cidade, das variáveis ambientais e da malha urbana." "@@TOKEN0@@" $false
Replace-Text "5840942 - Marco Aurélio Kondracki de Alcântara" "@@TOKEN1@@" $false
Replace-Text "Variável Ecológicano Ambiente Urbano; Enfoque Encômico e Impactos Ambientais." "@@TOKEN2@@" $false
Replace-Text "Elementos para estruturação ambiental da cidade. Variável ecológica no ambiente das atividades urbanas. A questão ambiental no urbanismo. A questão ambiental sob o enfoque econômico. Noções de higiene e saúde ambiental. A urbanização e os impactos ocasionados, principal enfoque da drenagem urbana." "@@TOKEN3@@" $false
Replace-Text "Aula expositiva e exercícios dirigidos." "@@TOKEN4@@" $false
Replace-Text "Média ponderada de exercícios e provas." "@@TOKEN5@@" $false
Replace-Text "Prova única com nota igual ou superior a 5,0." "@@TOKEN6@@" $false
Replace-Text "valle, C.R.*Drenagem Urbana. Ed. da Universidade e ABRH. 1995." "@@TOKEN7@@" $true

# ---- phase 2: replace each token with its destination text ----
Replace-Text "@@TOKEN0@@" "Variável Ecológicano Ambiente Urbano; Enfoque Encômico e Impactos Ambientais." $false
Replace-Text "@@TOKEN1@@" "Propiciar uma integração entre os elementos de estruturação da cidade, das variáveis ambientais e da malha urbana." $false
Replace-Text "@@TOKEN2@@" "Elementos para estruturação ambiental da cidade. Variável ecológica no ambiente das atividades urbanas. A questão ambiental no urbanismo. A questão ambiental sob o enfoque econômico. Noções de higiene e saúde ambiental. A urbanização e os impactos ocasionados, principal enfoque da drenagem urbana." $false
Replace-Text "@@TOKEN3@@" "Aula expositiva e exercícios dirigidos." $false
Replace-Text "@@TOKEN4@@" "Média ponderada de exercícios e provas." $false
Replace-Text "@@TOKEN5@@" "Prova única com nota igual ou superior a 5,0." $false
Replace-Text "@@TOKEN6@@" "valle, C.R. Qualidade ambiental: o desafio de ser competitivo protegendo o meio ambiente. Pioneira. 1995.^lDonaire, D.. Gestão ambiental na empresa. Atlas. 2a. edição. 1999.^lWinter, G.. Gestão e ambiente. Modelo prático de integração empresarial. Texto Editora, Lisboa. 1992.^lTucci, C.E., Porto, R.M., L.L. e Barros, M.T. org.. Drenagem Urbana. Ed. da Universidade e ABRH. 1995." $false
Replace-Text "@@TOKEN7@@" "5840942 - Marco Aurélio Kondracki de Alcântara" $false
